$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "NBOX MARATHON Over Ear Bluetooth Neckband 20 Hours Playback IPX5(Splash & Sweat Proof) Passive noise cancellation -Bluetooth Silver"

$ws.Range("A3").Value = "NBOX Buzz TWS On Ear True Wireless (TWS) 20 Hours Playback IPX5(Splash & Sweat Proof) Passive noise cancellation -Bluetooth Version 5.1 Black"
$ws.Range("B3").Value = "Rs. 749"

$ws.Range("A4").Value = "boAt Airdopes 121v2 On Ear True Wireless (TWS) 14 Hours Playback IPX7(Water Resistant) Active Noise cancellation -Bluetooth V 5.0 Black"
$ws.Range("B4").Value = "Rs. 1,299"

$ws.Range("A5").Value = "VEhop Power Bank Earbuds, On Ear True Wireless (TWS) 280 Hours Playback IPX4(Splash & Sweat Proof) Passive noise cancellation -Bluetooth V 5.2 Black"
$ws.Range("B5").Value = "Rs. 999"
